$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

$ws.Range("D2").Value = "66.488.32"
$ws.Range("E2").Value = "  +5.10%  "
$ws.Range("D3").Value = "3.849.22"
$ws.Range("E3").Value = "  +9.48%  "
Set-TextValue "D4" "0.999"
$ws.Range("E4").Value = "  -0.38%  "
Set-TextValue "D5" "427.39"
$ws.Range("E5").Value = "  +9.41%  "
$ws.Range("E6").Value = "  +8.94%  "
$ws.Range("D7").Value = "3.838.25"
$ws.Range("E7").Value = "  +9.43%  "
Set-TextValue "D8" "0.616"
$ws.Range("E8").Value = "  +5.09%  "
$ws.Range("E9").Value = "  -0.08%  "
Set-TextValue "D10" "0.733"
$ws.Range("E10").Value = "  +8.53%  "
Set-TextValue "D11" "0.159"
$ws.Range("E11").Value = "  +6.08%  "
$ws.Range("E12").Value = "  +2.13%  "
Set-TextValue "D13" "42.00"
$ws.Range("E13").Value = "  +8.60%  "
Set-TextValue "D14" "10.44"
$ws.Range("E14").Value = "  +13.86%  "
$ws.Range("D15").Value = "4.455.79"
$ws.Range("E15").Value = "  +9.65%  "
Set-TextValue "D16" "15.91"
$ws.Range("E16").Value = "  +25.52%  "
$ws.Range("D17").Value = "3.925.04"
$ws.Range("E17").Value = "  +11.94%  "
$ws.Range("E18").Value = "  +1.45%  "
Set-TextValue "D19" "20.12"
$ws.Range("E19").Value = "  +7.33%  "
$ws.Range("E20").Value = "  +8.45%  "
$ws.Range("D21").Value = "66.699.51"
$ws.Range("E21").Value = "  +5.30%  "
Set-TextValue "D22" "416.11"
$ws.Range("E22").Value = "  +5.82%  "
Set-TextValue "D23" "15.13"
$ws.Range("E23").Value = "  +9.42%  "
Set-TextValue "D24" "85.23"
$ws.Range("E24").Value = "  +5.64%  "
$ws.Range("E25").Value = "  +8.93%  "
Set-TextValue "D26" "37.58"
$ws.Range("E26").Value = "  +13.71%  "
Set-TextValue "D27" "10.07"
$ws.Range("E27").Value = "  +14.95%  "
Set-TextValue "D28" "3.28"
$ws.Range("E28").Value = "  +10.00%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D29" "9.34"
$ws.Range("E29").Value = "  +36.77%  "
$ws.Range("B30").Value = "LEO"
$ws.Range("C30").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D30" "5.35"
$ws.Range("E30").Value = "  +4.12%  "
Set-TextValue "D31" "719.74"
$ws.Range("E31").Value = "  +7.70%  "
Set-TextValue "D32" "13.83"
$ws.Range("E32").Value = "  +16.16%  "
$ws.Range("E33").Value = "  +14.78%  "
$ws.Range("E34").Value = "  +6.50%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  +44.93%  "
Set-TextValue "D37" "39.29"
$ws.Range("E37").Value = "  +7.58%  "
Set-TextValue "D38" "0.151"
$ws.Range("E38").Value = "  +1.45%  "
Set-TextValue "D39" "55.77"
$ws.Range("D40").Value = "0.0₃0752"
$ws.Range("E40").Value = "  +19.90%  "
Set-TextValue "D41" "0.0467"
$ws.Range("E41").Value = "  +7.23%  "
$ws.Range("E42").Value = "  +6.92%  "
$ws.Range("E43").Value = "  +0.27%  "
$ws.Range("E44").Value = "  +6.83%  "
$ws.Range("E45").Value = "  +4.74%  "
$ws.Range("E46").Value = "  +10.84%  "
Set-TextValue "D47" "0.322"
$ws.Range("E47").Value = "  +16.52%  "
$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D48" "2.64"
$ws.Range("E48").Value = "  +5.99%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D49" "2.87"
$ws.Range("E49").Value = "  +5.95%  "
$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D50" "2.06"
$ws.Range("E50").Value = "  +5.77%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D51" "142.78"
$ws.Range("E51").Value = "  +2.48%  "
